$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1912.4546
$ws.Range("I11").Value = 1912.4546
$ws.Range("K11").Value = 1912.4546
$ws.Range("M11").Value = -1772.4546
$ws.Range("H17").Value = 3131.6365
$ws.Range("J17").Value = 3298.3416
$ws.Range("L17").Value = 9895.024800000001
$ws.Range("N17").Value = -10231.0248
$ws.Range("H39").Value = 179.28
$ws.Range("I39").Value = 63.75
$ws.Range("J39").Value = 285.92307
$ws.Range("K39").Value = 191.25
$ws.Range("L39").Value = 857.7692099999999
$ws.Range("M39").Value = 104.75
$ws.Range("N39").Value = -1449.76921
$ws.Range("H70").Value = 1315.5333
$ws.Range("I70").Value = 976.4286
$ws.Range("J70").Value = 1612.25
$ws.Range("K70").Value = 2929.2858
$ws.Range("L70").Value = 4836.75
$ws.Range("M70").Value = -2659.2858
$ws.Range("N70").Value = -5376.75
$ws.Range("H73").Value = 1315.5333
$ws.Range("I73").Value = 976.4286
$ws.Range("J73").Value = 1612.25
$ws.Range("K73").Value = 2929.2858
$ws.Range("L73").Value = 4836.75
$ws.Range("M73").Value = -1993.2858
$ws.Range("N73").Value = -6708.75
$ws.Range("H92").Value = 803.29034
$ws.Range("I92").Value = 880.3913
$ws.Range("J92").Value = 581.625
$ws.Range("K92").Value = 880.3913
$ws.Range("L92").Value = 581.625
$ws.Range("M92").Value = 367.6087
$ws.Range("N92").Value = -3077.625
$ws.Range("H94").Value = 852.7143
$ws.Range("I94").Value = 852.7143
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 852.7143
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -401.7143
$ws.Range("N94").Value = ""
$ws.Range("H99").Value = 280
$ws.Range("J99").Value = 500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H103").Value = 506.24
$ws.Range("I103").Value = 531.7619
$ws.Range("J103").Value = 372.25
$ws.Range("K103").Value = 1595.2857
$ws.Range("L103").Value = 1116.75
$ws.Range("M103").Value = -1009.2857
$ws.Range("N103").Value = -2288.75
$ws.Range("H106").Value = 3417.111
$ws.Range("J106").Value = 4500
$ws.Range("L106").Value = 4500
$ws.Range("N106").Value = -5762
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3516.8784
$ws.Range("I32").Value = 2801.9849
$ws.Range("K32").Value = 2801.9849
$ws.Range("M32").Value = -2514.9849
$ws.Range("H45").Value = 7995373
$ws.Range("I45").Value = 14387166
$ws.Range("J45").Value = 5632.25
$ws.Range("K45").Value = 14387166
$ws.Range("L45").Value = 5632.25
$ws.Range("M45").Value = -14386789
$ws.Range("N45").Value = -6386.25
$ws.Range("H74").Value = 116200.95
$ws.Range("I74").Value = 38257.38
$ws.Range("J74").Value = 321688.53
$ws.Range("K74").Value = 38257.38
$ws.Range("L74").Value = 321688.53
$ws.Range("M74").Value = -37383.38
$ws.Range("N74").Value = -323436.53
$ws.Range("H77").Value = 116200.95
$ws.Range("I77").Value = 38257.38
$ws.Range("J77").Value = 321688.53
$ws.Range("K77").Value = 191286.9
$ws.Range("L77").Value = 1608442.65
$ws.Range("M77").Value = -186918.9
$ws.Range("N77").Value = -1617178.65
$ws.Range("H88").Value = 2636.7144
$ws.Range("I88").Value = 1290
$ws.Range("K88").Value = 1290
$ws.Range("M88").Value = -884
$ws.Range("H91").Value = 2636.7144
$ws.Range("I91").Value = 1290
$ws.Range("K91").Value = 1290
$ws.Range("M91").Value = 114
$ws.Range("H97").Value = 1545193.2
$ws.Range("I97").Value = 1802503.9
$ws.Range("J97").Value = 1330
$ws.Range("K97").Value = 1802503.9
$ws.Range("L97").Value = 1330
$ws.Range("M97").Value = -1802007.9
$ws.Range("N97").Value = -2322
$ws.Range("H122").Value = 402241.78
$ws.Range("I122").Value = 1345.6818
$ws.Range("J122").Value = 2607170.5
$ws.Range("K122").Value = 4037.0454
$ws.Range("L122").Value = 7821511.5
$ws.Range("M122").Value = -1587.0454
$ws.Range("N122").Value = -7826411.5
$ws.Range("H132").Value = 1966.2858
$ws.Range("I132").Value = 1412
$ws.Range("K132").Value = 4236
$ws.Range("M132").Value = -1706
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1339.0416
$ws.Range("I20").Value = 1039.8182
$ws.Range("K20").Value = 1039.8182
$ws.Range("M20").Value = -792.8181999999999
$ws.Range("H86").Value = 3454058.5
$ws.Range("I86").Value = 3577346.2
$ws.Range("K86").Value = 3577346.2
$ws.Range("M86").Value = -3576223.2
$ws.Range("H89").Value = 3454058.5
$ws.Range("I89").Value = 3577346.2
$ws.Range("K89").Value = 17886731
$ws.Range("M89").Value = -17881115
$ws.Range("H99").Value = 8931136
$ws.Range("I99").Value = 17859276
$ws.Range("K99").Value = 17859276
$ws.Range("M99").Value = -17857778
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4389.6665
$ws.Range("J99").Value = 5248.75
$ws.Range("L99").Value = 5248.75
$ws.Range("N99").Value = -8244.75
$ws.Range("H126").Value = 4389.6665
$ws.Range("J126").Value = 5248.75
$ws.Range("L126").Value = 15746.25
$ws.Range("N126").Value = -20686.25
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 35764.83
$ws.Range("I5").Value = 734.3125
$ws.Range("J5").Value = 78879.30499999999
$ws.Range("K5").Value = 2202.9375
$ws.Range("L5").Value = 236637.915
$ws.Range("M5").Value = -2090.9375
$ws.Range("N5").Value = -236861.915
$ws.Range("H86").Value = 174.42857
$ws.Range("I86").Value = 116
$ws.Range("J86").Value = 197.8
$ws.Range("K86").Value = 348
$ws.Range("L86").Value = 593.4000000000001
$ws.Range("M86").Value = 838
$ws.Range("N86").Value = -2965.4
$ws.Range("H89").Value = 174.42857
$ws.Range("I89").Value = 116
$ws.Range("J89").Value = 197.8
$ws.Range("K89").Value = 1044
$ws.Range("L89").Value = 1780.2
$ws.Range("M89").Value = 4884
$ws.Range("N89").Value = -13636.2
$ws.Range("H131").Value = 3806.2632
$ws.Range("J131").Value = 4544
$ws.Range("L131").Value = 13632
$ws.Range("N131").Value = -23712
$ws.Range("H135").Value = 35764.83
$ws.Range("I135").Value = 734.3125
$ws.Range("J135").Value = 78879.30499999999
$ws.Range("K135").Value = 6608.8125
$ws.Range("L135").Value = 709913.7449999999
$ws.Range("M135").Value = -4073.8125
$ws.Range("N135").Value = -714983.7449999999
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4070324.8
$ws.Range("I80").Value = 8132047.5
$ws.Range("J80").Value = 8602
$ws.Range("K80").Value = 8132047.5
$ws.Range("L80").Value = 8602
$ws.Range("M80").Value = -8131049.5
$ws.Range("N80").Value = -10598
$ws.Range("H83").Value = 4070324.8
$ws.Range("I83").Value = 8132047.5
$ws.Range("J83").Value = 8602
$ws.Range("K83").Value = 40660237.5
$ws.Range("L83").Value = 43010
$ws.Range("M83").Value = -40655245.5
$ws.Range("N83").Value = -52994
$ws.Range("H97").Value = 1191543.9
$ws.Range("I97").Value = 1323870.1
$ws.Range("J97").Value = 608
$ws.Range("K97").Value = 1323870.1
$ws.Range("L97").Value = 608
$ws.Range("M97").Value = -1323374.1
$ws.Range("N97").Value = -1600
$ws.Range("H102").Value = 10459407
$ws.Range("I102").Value = 37039664
$ws.Range("K102").Value = 37039664
$ws.Range("M102").Value = -37038042
# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3344.2083
$ws.Range("I7").Value = 1891.5625
$ws.Range("J7").Value = 6249.5
$ws.Range("K7").Value = 1891.5625
$ws.Range("L7").Value = 6249.5
$ws.Range("M7").Value = -1779.5625
$ws.Range("N7").Value = -6473.5
$ws.Range("H16").Value = 755.5238000000001
$ws.Range("I16").Value = 783.3
$ws.Range("K16").Value = 783.3
$ws.Range("M16").Value = -613.3
$ws.Range("H55").Value = 1563.0769
$ws.Range("I55").Value = 1301.7916
$ws.Range("K55").Value = 1301.7916
$ws.Range("M55").Value = -1128.7916
$ws.Range("H82").Value = 5051680.5
$ws.Range("J82").Value = 1322
$ws.Range("L82").Value = 1322
$ws.Range("N82").Value = -2044
$ws.Range("H85").Value = 5051680.5
$ws.Range("J85").Value = 1322
$ws.Range("L85").Value = 1322
$ws.Range("N85").Value = -3818
$ws.Range("H93").Value = 25643110
$ws.Range("J93").Value = 1411.75
$ws.Range("L93").Value = 1411.75
$ws.Range("N93").Value = -3907.75
$ws.Range("H126").Value = 3344.2083
$ws.Range("I126").Value = 1891.5625
$ws.Range("J126").Value = 6249.5
$ws.Range("K126").Value = 5674.6875
$ws.Range("L126").Value = 18748.5
$ws.Range("M126").Value = -3204.6875
$ws.Range("N126").Value = -23688.5
$ws.Range("H132").Value = 4549.1177
$ws.Range("I132").Value = 3751.9395
$ws.Range("J132").Value = 6010.6113
$ws.Range("K132").Value = 11255.8185
$ws.Range("L132").Value = 18031.8339
$ws.Range("M132").Value = -8725.818499999999
$ws.Range("N132").Value = -23091.8339
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8653.923000000001
$ws.Range("I62").Value = 1899.5
$ws.Range("K62").Value = 1899.5
$ws.Range("M62").Value = -1275.5
$ws.Range("H65").Value = 8653.923000000001
$ws.Range("I65").Value = 1899.5
$ws.Range("K65").Value = 9497.5
$ws.Range("M65").Value = -6377.5
$ws.Range("H81").Value = 11112007
$ws.Range("I81").Value = 11905579
$ws.Range("K81").Value = 23811158
$ws.Range("M81").Value = -23810097
$ws.Range("H84").Value = 11112007
$ws.Range("I84").Value = 11905579
$ws.Range("K84").Value = 119055790
$ws.Range("M84").Value = -119050486
$ws.Range("H96").Value = 1498.2858
$ws.Range("I96").Value = 1352
$ws.Range("J96").Value = 1693.3334
$ws.Range("K96").Value = 1352
$ws.Range("L96").Value = 1693.3334
$ws.Range("M96").Value = 21
$ws.Range("N96").Value = -4439.3334
$ws.Range("H126").Value = 1870.1765
$ws.Range("I126").Value = 2064.5
$ws.Range("K126").Value = 6193.5
$ws.Range("M126").Value = -3723.5
